$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, pushing existing rows down.
$ws.Rows.Item(2).Insert()

# The insert copies formatting from the row above (the bold header row) -
# clear it so the new data row matches the plain (unstyled) data rows.
$ws.Range("A2:B2").ClearFormats()

# Populate the newly inserted row 2 with the new commit entry.
$ws.Range("A2").Value = "Enhancement: Added GitHub commit fetch + integrated backend logic"
$ws.Range("B2").Value = "feature"
